$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 2 (F2:K2 "Hiver/Été/Année" sub-header row) is removed entirely;
# this shifts the four data rows (old 3-6) up to become rows 2-5.
$ws.Rows.Item(2).Delete()

# Clear any stale formatting on the header row before re-writing it, so the
# first five (new) columns end up with the default style.
$ws.Range("A1:K1").ClearFormats()

# New header row: idx, idx2, Name, Date Start, Date End, then the existing
# measurement headers (now re-labelled).
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give F1:K1 the same (9pt Arial) font used by the rest of the header/data
# cells, matching the new dedicated style.
$hdr = $ws.Range("F1:K1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 9

# Match the author's recorded selection after the edit (A2:K2, anchored A2).
$ws.Range("A2:K2").Select()
